$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "day 6" row (row 7) content is being moved down to become the
# "day 24" row (row 25), which previously only had the day number filled in.
$ws.Range("B25").Value = $ws.Range("B7").Text
$ws.Range("C25").Value = $ws.Range("C7").Text
$ws.Range("D25").Value = $ws.Range("D7").Text
$ws.Range("E25").Value = $ws.Range("E7").Text

# Row 7 ("day 6") now gets replaced with a brand new journal entry.
$ws.Range("B7").Value = "Ich zeige mich."
$ws.Range("C7").Value = "Ich habe Tanja's Aufruf im PM-Call vom 27.11 umgesetzt und aufgeschrieben, was ich erreicht habe und auf die Frage ""Wie will ich es haben?"" habe ich auch alles usäprätscht wie von Tanja gefordert. Heute zeige ich mich damit in der Gruppe. Sich so zu zeigen, löst bei mir die Frage aus, wen könnte das wirklich interessieren und ob ich mit dem erreichten jetzt bluffig rüberkomme. Das sind die Gefühle, die ich damit habe. Aber ich mache es jetzt trotzdem."
$ws.Range("D7").Value = "Mut"
$ws.Range("E7").Value = "https://ih1.redbubble.net/image.4999092252.0596/raf,360x360,075,t,fafafa:ca443f4786.jpg"

# Undo any automatic row-height override the engine applied while writing the
# long "text" cells, so the rows keep using the sheet's default row height
# exactly like the original file.
$ws.Rows.Item(7).EntireRow.AutoFit()
$ws.Rows.Item(25).EntireRow.AutoFit()

# Match the active selection shown in the saved workbook.
$ws.Range("C7").Select()
